$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price column D, Volume(1h) column E)
# D-column cells whose new text would otherwise be auto-recognized as a number
# are forced to remain text via NumberFormat "@" before assignment, matching the
# original inline-string (text) cell type.

$numericPriceRows = @(4, 5, 8, 9, 12, 14, 15, 16, 17, 19, 20, 22, 25, 28, 29, 30, 31, 33, 34, 35, 36, 37, 38, 39, 40, 41, 43, 44, 45, 46, 47, 48, 50, 51)
foreach ($r in $numericPriceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = '29.690.69'
$ws.Cells.Item(2, 5).Value = '  -2.93%  '
$ws.Cells.Item(3, 4).Value = '2.094.74'
$ws.Cells.Item(3, 5).Value = '  -2.21%  '
$ws.Cells.Item(4, 4).Value = '1.011'
$ws.Cells.Item(4, 5).Value = '  +0.18%  '
$ws.Cells.Item(5, 4).Value = '345.04'
$ws.Cells.Item(5, 5).Value = '  -1.89%  '
$ws.Cells.Item(6, 5).Value = '  +0.13%  '
$ws.Cells.Item(8, 4).Value = '0.4391'
$ws.Cells.Item(8, 5).Value = '  -3.96%  '
$ws.Cells.Item(9, 4).Value = '52.59'
$ws.Cells.Item(9, 5).Value = '  -2.26%  '
$ws.Cells.Item(10, 5).Value = '  +0.65%  '
$ws.Cells.Item(11, 5).Value = '  -1.09%  '
$ws.Cells.Item(12, 4).Value = '24.84'
$ws.Cells.Item(12, 5).Value = '  -2.18%  '
$ws.Cells.Item(13, 4).Value = '2.084.03'
$ws.Cells.Item(13, 5).Value = '  -2.67%  '
$ws.Cells.Item(14, 4).Value = '8.284'
$ws.Cells.Item(14, 5).Value = '  +1.54%  '
$ws.Cells.Item(15, 4).Value = '6.742'
$ws.Cells.Item(15, 5).Value = '  -2.04%  '
$ws.Cells.Item(16, 4).Value = '99.34'
$ws.Cells.Item(16, 5).Value = '  -2.62%  '
$ws.Cells.Item(17, 4).Value = '0.00001151'
$ws.Cells.Item(17, 5).Value = '  -1.83%  '
$ws.Cells.Item(18, 5).Value = '  +0.08%  '
$ws.Cells.Item(19, 4).Value = '20.87'
$ws.Cells.Item(19, 5).Value = '  +5.42%  '
$ws.Cells.Item(20, 4).Value = '0.06668'
$ws.Cells.Item(20, 5).Value = '  -0.58%  '
$ws.Cells.Item(21, 5).Value = '  +0.06%  '
$ws.Cells.Item(22, 4).Value = '6.190'
$ws.Cells.Item(22, 5).Value = '  -2.63%  '
$ws.Cells.Item(23, 4).Value = '29.735.13'
$ws.Cells.Item(23, 5).Value = '  -3.14%  '
$ws.Cells.Item(24, 5).Value = '  -2.69%  '
$ws.Cells.Item(25, 4).Value = '2.322'
$ws.Cells.Item(25, 5).Value = '  -2.55%  '
$ws.Cells.Item(26, 4).Value = '2.333.51'
$ws.Cells.Item(26, 5).Value = '  -1.60%  '
$ws.Cells.Item(27, 5).Value = '  -3.09%  '
$ws.Cells.Item(28, 4).Value = '2.523'
$ws.Cells.Item(28, 5).Value = '  -4.64%  '
$ws.Cells.Item(29, 4).Value = '161.96'
$ws.Cells.Item(29, 5).Value = '  -1.64%  '
$ws.Cells.Item(30, 4).Value = '133.16'
$ws.Cells.Item(30, 5).Value = '  -3.05%  '
$ws.Cells.Item(31, 4).Value = '1.131'
$ws.Cells.Item(31, 5).Value = '  -7.26%  '
$ws.Cells.Item(32, 5).Value = '  -2.90%  '
$ws.Cells.Item(33, 4).Value = '1.649'
$ws.Cells.Item(33, 5).Value = '  -1.38%  '
$ws.Cells.Item(34, 4).Value = '6.167'
$ws.Cells.Item(34, 5).Value = '  -3.32%  '
$ws.Cells.Item(35, 4).Value = '3.936'
$ws.Cells.Item(35, 5).Value = '  -2.13%  '
$ws.Cells.Item(36, 4).Value = '6.191'
$ws.Cells.Item(36, 5).Value = '  +0.36%  '
$ws.Cells.Item(37, 4).Value = '10.26'
$ws.Cells.Item(37, 5).Value = '  -1.78%  '
$ws.Cells.Item(38, 4).Value = '0.02570'
$ws.Cells.Item(38, 5).Value = '  -3.13%  '
$ws.Cells.Item(39, 4).Value = '0.06695'
$ws.Cells.Item(39, 5).Value = '  -4.50%  '
$ws.Cells.Item(40, 4).Value = '12.44'
$ws.Cells.Item(40, 5).Value = '  -1.94%  '
$ws.Cells.Item(41, 4).Value = '0.6853'
$ws.Cells.Item(41, 5).Value = '  -2.26%  '
$ws.Cells.Item(42, 5).Value = '  -4.88%  '
$ws.Cells.Item(43, 4).Value = '1.301'
$ws.Cells.Item(43, 5).Value = '  +1.76%  '
$ws.Cells.Item(44, 4).Value = '0.6640'
$ws.Cells.Item(44, 5).Value = '  +2.55%  '
$ws.Cells.Item(45, 4).Value = '14.28'
$ws.Cells.Item(45, 5).Value = '  -3.56%  '
$ws.Cells.Item(46, 4).Value = '2.315'
$ws.Cells.Item(46, 5).Value = '  -2.06%  '
$ws.Cells.Item(47, 4).Value = '3.635'
$ws.Cells.Item(47, 5).Value = '  -3.06%  '
$ws.Cells.Item(48, 4).Value = '0.00000000355'
$ws.Cells.Item(48, 5).Value = '  -3.58%  '
$ws.Cells.Item(49, 5).Value = '  -2.84%  '
$ws.Cells.Item(50, 4).Value = '82.12'
$ws.Cells.Item(50, 5).Value = '  -1.29%  '
$ws.Cells.Item(51, 4).Value = '0.3295'
$ws.Cells.Item(51, 5).Value = '  +0.32%  '
